$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, shifting existing rows 21-86 down to 22-87.
$ws.Rows("21:21").Insert()

# Populate the newly inserted row 21 with the new record's data.
$ws.Range("A21").Value = 5
$ws.Range("B21").Value = "Macroferia Regional de Talca"
$ws.Range("C21").Value = "Maule"
$ws.Range("D21").Value = 44592
$ws.Range("E21").Value = 7
$ws.Range("F21").Value = 100112022
$ws.Range("G21").Value = "Arveja Verde"
$ws.Range("H21").Value = "Perfection"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 250
$ws.Range("K21").Value = 25000
$ws.Range("L21").Value = 25000
$ws.Range("M21").Value = 25000
$ws.Range("N21").Value = "$/saco 25 kilos"
$ws.Range("O21").Value = "Carahue"
$ws.Range("P21").Value = 1000
$ws.Range("Q21").Value = 25
$ws.Range("R21").Value = "Hortaliza"
